$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# New progress as of 04-Nov-2025:
#  - PERIOD TO EXPIRE (column H) for rows 3-8 each drop by one day
#  - LAST UPDATE (column I) for rows 3-8 moves from 03-Nov-2025 to 04-Nov-2025
#    (kept as literal text via a leading apostrophe so it stays a text value,
#    matching the existing "LAST UPDATE" column rather than becoming a date serial)

for ($row = 3; $row -le 8; $row++) {
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)
    $iCell.Formula = "'04-Nov-2025"
}
